# Add SFX rows for the new "ItemSystem" team / "Item" class, then re-sort
# the SFX sheet by the generated FileName (column G), matching the sheet's
# existing sortState.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SFX")

# New rows describing the Item System sounds.
$newRows = @(
    @("ItemSystem", "Item", "Bomb",         "GreatShoot.wav"),
    @("ItemSystem", "Item", "Invincible_10s","Zhoararang.wav"),
    @("ItemSystem", "Item", "ItemGet1",      "AscendingScales1.wav"),
    @("ItemSystem", "Item", "SpeedUp_10s",   "ShipTakeOff.wav"),
    @("ItemSystem", "Item", "SubShip_10s",   "ShipLandingAndTakeOff.wav")
)

$startRow = 14
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = "S"
    $ws.Cells.Item($r, 6).Value = "O"
    $ws.Cells.Item($r, 7).Formula = '=_xlfn.CONCAT("SFX_",A' + $r + ',"_",B' + $r + ',"_",C' + $r + ',"_",D' + $r + ')'
}

# Re-sort the full data range (rows 2-18) ascending by column G, mirroring
# the workbook's existing sortState (ref="A2:G20" sortCondition ref="G2:G20").
$sortRange = $ws.Range("A2:G18")
$sortKey = $ws.Range("G2:G18")
$sortRange.Sort($sortKey)

$ws.Range("A1:G18").Columns.AutoFit()
